# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" data table: a handful of countries report new
# totals, which re-sorts a few neighbouring rows (the sheet is kept sorted
# descending by "Casos totales" / column B), and the "last updated" footer
# timestamp moves from 20:22 to 20:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- footer timestamp (A1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Abril de 2020 a las 20:52"

# --- Estados Unidos (row 4): Muertes hoy ---
$ws.Cells.Item(4,6).Value = 8876

# --- Brasil (row 19) ---
$ws.Cells.Item(19,2).Value = 11670
$ws.Cells.Item(19,3).Value = 416
$ws.Cells.Item(19,5).Value = 11030
$ws.Cells.Item(19,7).Value = 27
$ws.Cells.Item(19,8).Value = 513

# --- Israel (row 21) ---
$ws.Cells.Item(21,2).Value = 8904
$ws.Cells.Item(21,3).Value = 474
$ws.Cells.Item(21,5).Value = 8262
$ws.Cells.Item(21,6).Value = 140
$ws.Cells.Item(21,7).Value = 8
$ws.Cells.Item(21,8).Value = 57

# --- Kazajistan overtakes Camerun (rows 75-76) ---
$ws.Cells.Item(75,1).Value = "Kazajistan"
$ws.Cells.Item(75,2).Value = 662
$ws.Cells.Item(75,3).Value = 78
$ws.Cells.Item(75,4).Value = 46
$ws.Cells.Item(75,5).Value = 610
$ws.Cells.Item(75,6).Value = 16
$ws.Cells.Item(75,7).Value = 0
$ws.Cells.Item(75,8).Value = 6

$ws.Cells.Item(76,1).Value = "Camerun"
$ws.Cells.Item(76,2).Value = 658
$ws.Cells.Item(76,3).Value = 8
$ws.Cells.Item(76,4).Value = 17
$ws.Cells.Item(76,5).Value = 632
$ws.Cells.Item(76,6).Value = 0
$ws.Cells.Item(76,7).Value = 0
$ws.Cells.Item(76,8).Value = 9

# --- Uzbekistan jumps ahead of Costa Rica and Uruguay (rows 86-88) ---
$ws.Cells.Item(86,1).Value = "Uzbekistan"
$ws.Cells.Item(86,2).Value = 457
$ws.Cells.Item(86,3).Value = 115
$ws.Cells.Item(86,4).Value = 30
$ws.Cells.Item(86,5).Value = 425
$ws.Cells.Item(86,6).Value = 8
$ws.Cells.Item(86,7).Value = 0
$ws.Cells.Item(86,8).Value = 2

$ws.Cells.Item(87,1).Value = "Costa Rica"
$ws.Cells.Item(87,2).Value = 454
$ws.Cells.Item(87,3).Value = 0
$ws.Cells.Item(87,4).Value = 16
$ws.Cells.Item(87,5).Value = 436
$ws.Cells.Item(87,6).Value = 14
$ws.Cells.Item(87,7).Value = 0
$ws.Cells.Item(87,8).Value = 2

$ws.Cells.Item(88,1).Value = "Uruguay"
$ws.Cells.Item(88,2).Value = 406
$ws.Cells.Item(88,3).Value = 0
$ws.Cells.Item(88,4).Value = 104
$ws.Cells.Item(88,5).Value = 296
$ws.Cells.Item(88,6).Value = 14
$ws.Cells.Item(88,7).Value = 0
$ws.Cells.Item(88,8).Value = 6

# --- Ruanda overtakes Trinidad yTobago (rows 127-128) ---
$ws.Cells.Item(127,1).Value = "Ruanda"
$ws.Cells.Item(127,2).Value = 105
$ws.Cells.Item(127,3).Value = 1
$ws.Cells.Item(127,4).Value = 4
$ws.Cells.Item(127,5).Value = 101
$ws.Cells.Item(127,6).Value = 0
$ws.Cells.Item(127,7).Value = 0
$ws.Cells.Item(127,8).Value = 0

$ws.Cells.Item(128,1).Value = "Trinidad yTobago"
$ws.Cells.Item(128,2).Value = 105
$ws.Cells.Item(128,3).Value = 1
$ws.Cells.Item(128,4).Value = 1
$ws.Cells.Item(128,5).Value = 96
$ws.Cells.Item(128,6).Value = 0
$ws.Cells.Item(128,7).Value = 1
$ws.Cells.Item(128,8).Value = 8

# --- Guinea Ecuatorial (row 167) ---
$ws.Cells.Item(167,4).Value = 3
$ws.Cells.Item(167,5).Value = 13

# --- Angola jumps ahead of Antigua y Barbuda, Mongolia, Dominica, Fiyi,
#     Santa Lucia (rows 169-174) ---
$ws.Cells.Item(169,1).Value = "Angola"
$ws.Cells.Item(169,2).Value = 16
$ws.Cells.Item(169,3).Value = 2
$ws.Cells.Item(169,4).Value = 2
$ws.Cells.Item(169,5).Value = 12
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 2

$ws.Cells.Item(170,1).Value = "Antigua y Barbuda"
$ws.Cells.Item(170,2).Value = 15
$ws.Cells.Item(170,3).Value = 0
$ws.Cells.Item(170,4).Value = 0
$ws.Cells.Item(170,5).Value = 15
$ws.Cells.Item(170,6).Value = 1
$ws.Cells.Item(170,7).Value = 0
$ws.Cells.Item(170,8).Value = 0

$ws.Cells.Item(171,1).Value = "Mongolia"
$ws.Cells.Item(171,2).Value = 15
$ws.Cells.Item(171,3).Value = 1
$ws.Cells.Item(171,4).Value = 2
$ws.Cells.Item(171,5).Value = 13
$ws.Cells.Item(171,6).Value = 0
$ws.Cells.Item(171,7).Value = 0
$ws.Cells.Item(171,8).Value = 0

$ws.Cells.Item(172,1).Value = "Dominica"
$ws.Cells.Item(172,2).Value = 14
$ws.Cells.Item(172,3).Value = 0
$ws.Cells.Item(172,4).Value = 0
$ws.Cells.Item(172,5).Value = 14
$ws.Cells.Item(172,6).Value = 0
$ws.Cells.Item(172,7).Value = 0
$ws.Cells.Item(172,8).Value = 0

$ws.Cells.Item(173,1).Value = "Fiyi"
$ws.Cells.Item(173,2).Value = 14
$ws.Cells.Item(173,3).Value = 2
$ws.Cells.Item(173,4).Value = 0
$ws.Cells.Item(173,5).Value = 14
$ws.Cells.Item(173,6).Value = 0
$ws.Cells.Item(173,7).Value = 0
$ws.Cells.Item(173,8).Value = 0

$ws.Cells.Item(174,1).Value = "Santa Lucia"
$ws.Cells.Item(174,2).Value = 14
$ws.Cells.Item(174,3).Value = 0
$ws.Cells.Item(174,4).Value = 1
$ws.Cells.Item(174,5).Value = 13
$ws.Cells.Item(174,6).Value = 0
$ws.Cells.Item(174,7).Value = 0
$ws.Cells.Item(174,8).Value = 0

# --- Liberia (row 175) ---
$ws.Cells.Item(175,2).Value = 14
$ws.Cells.Item(175,3).Value = 1
$ws.Cells.Item(175,5).Value = 8
